# Update row 8 (year 2025) metrics in metricas_recorrencia_anual
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1436
$ws.Range("D8").Value = 219
$ws.Range("E8").Value = 1217
$ws.Range("F8").Value = 8.982772764561116
$ws.Range("G8").Value = 84.74930362116991
$ws.Range("H8").Value = 15.25069637883008
